# Upload data: 15 January 2022
# Adds the December 2021 (date serial 44531 = 2021-12-01) unemployment rows
# to both the "Canada" sheet and the "Province" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Canada" — add row 25 (Canada, Dec-2021)
# ---------------------------------------------------------------------
$wsCanada = $wb.Worksheets.Item("Canada")

$wsCanada.Range("A25").Value = 44531
$wsCanada.Range("A25").NumberFormat = "d-mmm-yy"

$wsCanada.Range("B25").Value = "Canada"
$wsCanada.Range("B25").NumberFormat = "d-mmm-yy"

$wsCanada.Range("D25").Value = 1212.3
$wsCanada.Range("E25").Value = 1148.7
$wsCanada.Range("C25").Formula = "=(D25-E25)/E25*100"

# ---------------------------------------------------------------------
# Sheet "Province" — add rows 232-241 (each province, Dec-2021)
# ---------------------------------------------------------------------
$wsProvince = $wb.Worksheets.Item("Province")

$provinceRows = @(
    @{ Row = 232; Name = "Newfoundland & Labrador"; D = 29.5;   E = 31.4 },
    @{ Row = 233; Name = "Prince Edward Island";     D = 7;     E = 6.8 },
    @{ Row = 234; Name = "Nova Scotia";              D = 40.1;  E = 40.3 },
    @{ Row = 235; Name = "New Brunswick";            D = 31.4;  E = 29.7 },
    @{ Row = 236; Name = "Quebec";                   D = 211;   E = 236.7 },
    @{ Row = 237; Name = "Ontario";                  D = 491.9; E = 419.7 },
    @{ Row = 238; Name = "Manitoba";                 D = 36.4;  E = 35.1 },
    @{ Row = 239; Name = "Saskatchewan";              D = 32.6;  E = 36.3 },
    @{ Row = 240; Name = "Alberta";                  D = 181.4; E = 178.5 },
    @{ Row = 241; Name = "British Columbia";         D = 151;   E = 134.1 }
)

foreach ($r in $provinceRows) {
    $row = $r.Row

    $wsProvince.Range("A$row").Value = 44531
    $wsProvince.Range("A$row").NumberFormat = "d-mmm-yy"

    $wsProvince.Range("B$row").Value = $r.Name
    if ($row -eq 232) {
        # First province of the date block carries the date-style (matches
        # the source workbook, which styles the "Newfoundland & Labrador"
        # cell the same as the date column).
        $wsProvince.Range("B$row").NumberFormat = "d-mmm-yy"
    }

    $wsProvince.Range("D$row").Value = $r.D
    $wsProvince.Range("E$row").Value = $r.E
    $wsProvince.Range("C$row").Formula = "=(D$row-E$row)/E$row*100"
}

# ---------------------------------------------------------------------
# View state — mirror the selection/scroll position left by the author
# after appending the new rows.
# ---------------------------------------------------------------------
$wsCanada.Activate()
$wsCanada.Range("D14").Select()

$wsProvince.Activate()
$wsProvince.Range("D242").Select()

Write-Output "done"
